$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 1
$ws.Range("AB2").Value = 4
$ws.Range("AC2").Value = 8

# Row 4
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 1
$ws.Range("AC4").Value = 8

# Row 5
$ws.Range("X5").Value = 1
$ws.Range("Y5").Value = 1
$ws.Range("Z5").Value = 1
$ws.Range("AA5").Value = 1
$ws.Range("AC5").Value = 8

# Row 6
# Assign the classic "text" apostrophe prefix so the cell becomes an
# explicit empty-text cell (matching the original empty inlineStr cell)
# instead of being cleared back to a blank/empty cell.
$ws.Range("Z6").Value = "'"
$ws.Range("Z6").Style = "Normal"
$ws.Range("AA6").Value = "'"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").Value = 1
$ws.Range("AC6").Value = 2

# Row 7
$ws.Range("V7").Value = 1
$ws.Range("W7").Value = 1
$ws.Range("AC7").Value = 8
$ws.Range("AD7").Value = "Паритет"
